$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 285 (shifts existing rows 285-312 down to 288-315)
$ws.Rows.Item(285).Resize(3).Insert()

# Row 285 - Especial
$ws.Range("A285").Value = 9
$ws.Range("B285").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C285").Value = "Metropolitana"
$ws.Range("D285").Value = 45154
$ws.Range("D285").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E285").Value = 13
$ws.Range("F285").Value = "Fruta"
$ws.Range("G285").Value = 100107
$ws.Range("H285").Value = "Otros"
$ws.Range("I285").Value = 100107002
$ws.Range("J285").Value = "Chirimoya"
$ws.Range("K285").Value = "Cultivar IV Región"
$ws.Range("L285").Value = "Especial"
$ws.Range("M285").Value = 120
$ws.Range("N285").Value = 21600
$ws.Range("O285").Value = 21600
$ws.Range("P285").Value = 21600
$ws.Range("Q285").Value = "`$/bandeja 8 kilos"
$ws.Range("R285").Value = "Provincia del Elquí"
$ws.Range("S285").Value = 2700
$ws.Range("T285").Value = 8

# Row 286 - Primera
$ws.Range("A286").Value = 9
$ws.Range("B286").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C286").Value = "Metropolitana"
$ws.Range("D286").Value = 45154
$ws.Range("D286").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E286").Value = 13
$ws.Range("F286").Value = "Fruta"
$ws.Range("G286").Value = 100107
$ws.Range("H286").Value = "Otros"
$ws.Range("I286").Value = 100107002
$ws.Range("J286").Value = "Chirimoya"
$ws.Range("K286").Value = "Cultivar IV Región"
$ws.Range("L286").Value = "Primera"
$ws.Range("M286").Value = 160
$ws.Range("N286").Value = 17600
$ws.Range("O286").Value = 17600
$ws.Range("P286").Value = 17600
$ws.Range("Q286").Value = "`$/bandeja 8 kilos"
$ws.Range("R286").Value = "Provincia del Elquí"
$ws.Range("S286").Value = 2200
$ws.Range("T286").Value = 8

# Row 287 - Segunda
$ws.Range("A287").Value = 9
$ws.Range("B287").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C287").Value = "Metropolitana"
$ws.Range("D287").Value = 45154
$ws.Range("D287").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E287").Value = 13
$ws.Range("F287").Value = "Fruta"
$ws.Range("G287").Value = 100107
$ws.Range("H287").Value = "Otros"
$ws.Range("I287").Value = 100107002
$ws.Range("J287").Value = "Chirimoya"
$ws.Range("K287").Value = "Cultivar IV Región"
$ws.Range("L287").Value = "Segunda"
$ws.Range("M287").Value = 150
$ws.Range("N287").Value = 14400
$ws.Range("O287").Value = 14400
$ws.Range("P287").Value = 14400
$ws.Range("Q287").Value = "`$/bandeja 8 kilos"
$ws.Range("R287").Value = "Provincia del Elquí"
$ws.Range("S287").Value = 1800
$ws.Range("T287").Value = 8
